# Implements: doc(mysql): implements full join use alternative method
#
# The IDEA/PyCharm shortcut-reference sheet gains four new shortcut rows:
#   - ctrl+shift+j / "Join lines, 将两行合成一行"      (inserted right after ctrl+shift+enter)
#   - ctrl+enter    / "split line，将一行拆分成两行"    (inserted right after ctrl+shift+j)
#   - alt+shift+,   / "缩小所有编辑器的字体大小"         (inserted right after "将该行代码下移")
#   - alt+shift+.   / "增大所有编辑器的字体大小"         (inserted right after alt+shift+,)
# and the description for ctrl+shift+enter is reworded/expanded.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Reword the ctrl+shift+enter description (row 64, column B) and grow the
#    row so the extra wrapped line is visible.
$newDesc = "完成当行代码，如括号、分号可自动补全，补全后再次敲击，可跳下一行（如是空行）（不同于shift+enter新起一行)" + [char]10 + "Choose lookup item and invoke complete statement"
$ws.Range("B64").Value2 = $newDesc
$ws.Rows(64).RowHeight = 67.5

# 2) Insert two new rows right after row 64 for the "join lines" / "split
#    line" shortcuts. Inserting copies the formatting (fill/banding) of the
#    row above automatically, matching the existing striped table look.
$ws.Rows("65:66").Insert()
$ws.Range("A65").Value2 = "ctrl+shift+j"
$ws.Range("B65").Value2 = "Join lines, 将两行合成一行"
$ws.Range("A66").Value2 = "ctrl+enter"
$ws.Range("B66").Value2 = "split line，将一行拆分成两行"

# 3) Insert two more rows after "alt+shift+down / 将该行代码下移" (originally
#    row 72, now row 74 after the previous insert) for the editor
#    font-size shortcuts.
$ws.Rows("75:76").Insert()
$ws.Range("A75").Value2 = "alt+shift+,"
$ws.Range("B75").Value2 = "缩小所有编辑器的字体大小"
$ws.Range("A76").Value2 = "alt+shift+."
$ws.Range("B76").Value2 = "增大所有编辑器的字体大小"

# 4) Update the frozen-pane selection to match where editing left off.
$ws.Range("B73").Select()
